$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-7 from 45174 to 45175
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45175
}
